$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update header title / labels (shared string shifts are implicit; we just set text)
$ws.Range("D2").Value = "ESTADO DE CUENTA"
$ws.Range("B7").Value = "RAZON SOCIAL:"
$ws.Range("B11").Value = "VALOR MORA"
$ws.Range("E11").Value = 157509
$ws.Range("B13").Value = "Cant. Trabajadores"
$ws.Range("E13").Value = "Cant. Periodos"
$ws.Range("F13").Value = 3

# Swap the Novedad de Retiro / Novedad de Ingreso header columns
$ws.Range("H15").Value = "Novedad de Ingreso"
$ws.Range("I15").Value = "Novedad de Retiro"

# Update G17 value
$ws.Range("G17").Value = 2000000

# Insert a new row 18 (copy formatting of row 17), shifting everything below down
$ws.Rows("18:18").Insert()

$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "30290872"
$ws.Range("D18").Value = "OLGA LUCIA VALENCIA LOPEZ"
$ws.Range("E18").Value = "2508"
$ws.Range("F18").Value = 80000
$ws.Range("G18").Value = 2000000

$wb.Save()
